# Turn File Contents.xlsx - Sheet1 updates
#
# Started Orders.ToXml() and supporting methods. These are disabled in
# OrderWriter to use the old binary format until they are completed and
# tested.
#
# This updates the "Proposed" (right-hand) side of the documentation table
# on Sheet1: renames a couple of fields to match the new Orders/Intel
# classes, adds an explanatory note to RaceData, and fills in the three
# new "All*" columns that were previously blank for rows 4-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 2 - ConsoleState proposed layout: NewTurn renamed to Intel,
# RaceTurn renamed to Orders (Universe stays the same).
$ws.Range("J2").Value = "Intel"
$ws.Range("L2").Value = "Orders"

# Rows 4-6 previously had no value in column I (the "All*" class column);
# fill them in to mirror column H.
$ws.Range("I4").Value = "AllDesigns"
$ws.Range("I5").Value = "AllFleets"
$ws.Range("I6").Value = "AllMinefields"

# Row 9 - "AllStars (position only)" simplified to just "AllStars".
$ws.Range("I9").Value = "AllStars"

# Row 39 - clarify that Race.RaceData is not actually a RaceData object.
$ws.Range("D39").Value = "RaceData (not a RaceData object)"

# Reflect the cell that was being edited/viewed when the workbook was saved.
$ws.Range("D40").Select()
